$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TODAY")

# Update K and L columns for rows 2, 4, 5, 6, 7 with new values,
# and the totals row 8, per the commit's revised figures.
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 7.17

$ws.Range("K4").Value = 10
$ws.Range("L4").Value = 16.34

$ws.Range("K5").Value = 16.765
$ws.Range("L5").Value = 24.48

$ws.Range("K6").Value = 23
$ws.Range("L6").Value = 60.77

$ws.Range("K7").Value = 83.44799999999999
$ws.Range("L7").Value = 63.35

$ws.Range("K8").Value = 138.213
$ws.Range("L8").Value = 175.95
